$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text content changes on existing rows -------------------------------

# Row 10: "column I" -> "column H" (the proof table now lives in column H)
$ws.Range("A10").Value = "Respond with the response in column H below."

# B14 is the current challenge word; it changes to a new word.
$ws.Range("B14").Value = "FSMDA"

# --- Style-only tweaks (apply before copying formats to the new block so
# --- the copy below picks up the corrected formatting). -------------------

# E13 / E15 become right-aligned (matches F13/F15 etc.)
$ws.Range("E13").HorizontalAlignment = -4152   # xlRight
$ws.Range("E15").HorizontalAlignment = -4152   # xlRight

# --- New explanatory line for the fixed-clock scenario --------------------

$ws.Range("A22").Value = "Run Adventure with -d1-JAN-77 -t0:00 to lock the date and time and simplify the calculation."

# --- Build the second (fixed clock) proof table at rows 23-30 by cloning
# --- the formatting of the live table at rows 13-20, then filling in the
# --- static-date version of the formulas/values. ---------------------------

# Columns that are populated on every row of the live table (13-20): copy
# their formatting straight down to the new block (23-30).
foreach ($col in @("A", "B", "E", "F", "H")) {
    $ws.Range("${col}13:${col}20").Copy()
    $ws.Range("${col}23").PasteSpecial(-4122)   # xlPasteFormats
}
# Columns C, D, G are only populated (and only styled) on the header row
# (15); the data rows (16-20) are unstyled, which is also what rows 26-30
# should look like, so only the header row's formatting needs copying.
foreach ($col in @("C", "D", "G")) {
    $ws.Range("${col}15").Copy()
    $ws.Range("${col}25").PasteSpecial(-4122)   # xlPasteFormats
}
$excel.CutCopyMode = $false

# Row 23: header row (labels)
$ws.Range("A23").Value = "Time"
$ws.Range("B23").Value = "Challenge"
$ws.Range("E23").Value = "MAGNUM"
$ws.Range("F23").Value = "Minutes"
$ws.Range("H23").Value = "Response"

# Row 24: fixed clock of 1977-01-01 00:00 (serial 28126) instead of NOW()
$ws.Range("A24").Value = 28126
$ws.Range("A24").NumberFormat = "m/d/yy h:mm"
$ws.Range("B24").Value = "FSMDA"
$ws.Range("E24").Value = 11111
$ws.Range("F24").Formula = "=60*HOUR(A24)+MINUTE(A24)"
$ws.Range("H24").Formula = "=CONCAT(H26:H30)"

# Row 25: column headers for the letter-by-letter breakdown
$ws.Range("A25").Value = "Y"
$ws.Range("B25").Value = "letter"
$ws.Range("C25").Value = "VAL(Y)"
$ws.Range("D25").Value = "|delta|"
$ws.Range("E25").Value = "D"
$ws.Range("F25").Value = "T"
$ws.Range("G25").Value = "X"
$ws.Range("H25").Value = "letter"

# Rows 26-30: per-letter computation, mirroring rows 16-20 but sourced from
# row 24 (the fixed clock) instead of row 14 (NOW()).
$ws.Range("A26").Value = 1
$ws.Range("A27").Value = 2
$ws.Range("A28").Value = 3
$ws.Range("A29").Value = 4
$ws.Range("A30").Value = 5

for ($i = 0; $i -lt 5; $i++) {
    $r = 26 + $i
    $prev = $r - 1
    $next = $r + 1
    if ($i -eq 4) { $next = 26 }

    $ws.Range("B$r").Formula = "=UPPER(MID(B24, A$r, 1))"
    $ws.Range("C$r").Formula = "=CODE(UPPER(B$r)) - CODE(""A"") + 1"
    $ws.Range("D$r").Formula = "=ABS(C$r-C$next)"
    $ws.Range("G$r").Formula = "=MOD(D$r*MOD(E$r, 10)+MOD(F$r, 10), 26)+1"
    $ws.Range("H$r").Formula = "=CHAR(64+G$r)"

    if ($i -eq 0) {
        $ws.Range("E$r").Formula = "=E24"
        $ws.Range("F$r").Formula = "=40*FLOOR(F24/60, 1)+10*FLOOR(F24/10,1)"
    } else {
        $ws.Range("E$r").Formula = "=FLOOR(E$prev/10, 1)"
        $ws.Range("F$r").Formula = "=FLOOR(F$prev/10, 1)"
    }
}

# --- Selection moves to C14, matching the recorded UI state ---------------
$ws.Range("C14").Select()
